$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Participação no mercado do chai na América Latina" "Plano de Promoção do Chá Chai na América Latina"
Replace-Text "O Mystic Spice Premium Chai Tea é um tipo de chá aromático que teve origem na Índia e se tornou popular em todo o mundo." "O chá chai é uma bebida de chá temperada que se originou na Índia e se tornou popular em todo o mundo."
Replace-Text "O plano promocional e estratégia para o chai na América Latina tem como objetivos:" "O plano de promoção do chá Chai na América Latina visa atingir os seguintes objetivos:"
Replace-Text "O plano promocional e estratégia para o chai na América Latina utilizará uma combinação de táticas, incluindo:" "O plano de promoção do chá Chai na América Latina usará uma combinação de táticas, tais como:"
Replace-Text "O plano promocional e estratégia para o chai na América Latina será implementado ao longo de um período de 12 meses, com um orçamento de US`$100.000." "O plano de promoção do chá Chai na América Latina será implementado durante um período de 12 meses, com um orçamento de US `$ 100.000."
